$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '29.541.55'
Set-TextValue "E2" '  +1.01%  '

# Row 3
Set-TextValue "D3" '1.880.62'
Set-TextValue "E3" '  +1.04%  '

# Row 4
Set-TextValue "D4" '0.9988'
Set-TextValue "E4" '  -0.24%  '

# Row 5
Set-TextValue "D5" '0.7239'
Set-TextValue "E5" '  +3.09%  '

# Row 6
Set-TextValue "D6" '239.28'
Set-TextValue "E6" '  +0.66%  '

# Row 7
Set-TextValue "D7" '0.9986'
Set-TextValue "E7" '  -0.22%  '

# Row 8
Set-TextValue "D8" '0.07887'
Set-TextValue "E8" '  -2.74%  '

# Row 9
Set-TextValue "D9" '0.3087'
Set-TextValue "E9" '  +1.99%  '

# Row 10
Set-TextValue "D10" '25.37'
Set-TextValue "E10" '  +9.41%  '

# Row 11
Set-TextValue "D11" '0.08212'
Set-TextValue "E11" '  +0.69%  '

# Row 12
Set-TextValue "D12" '1.863.94'
Set-TextValue "E12" '  +1.43%  '

# Row 13
Set-TextValue "D13" '5.272'
Set-TextValue "E13" '  +2.21%  '

# Row 14
Set-TextValue "D14" '0.7266'
Set-TextValue "E14" '  +3.07%  '

# Row 15
Set-TextValue "D15" '89.65'
Set-TextValue "E15" '  +0.75%  '

# Row 16
Set-TextValue "D16" '29.479.51'
Set-TextValue "E16" '  +0.77%  '

# Row 17
Set-TextValue "D17" '5.853'
Set-TextValue "E17" '  +1.54%  '

# Row 18
Set-TextValue "D18" '0.000007875'
Set-TextValue "E18" '  +0.54%  '

# Row 19
Set-TextValue "B19" 'BitcoinCash'
Set-TextValue "C19" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue "D19" '242.12'
Set-TextValue "E19" '  +3.02%  '

# Row 20
Set-TextValue "B20" 'Avalanche'
Set-TextValue "C20" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D20" '13.38'
Set-TextValue "E20" '  +0.66%  '

# Row 21
Set-TextValue "B21" 'Dai'
Set-TextValue "C21" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D21" '0.9993'
Set-TextValue "E21" '  -0.12%  '

# Row 22
Set-TextValue "B22" 'WrappedliquidstakedEther2.0'
Set-TextValue "C22" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue "D22" '2.112.49'
Set-TextValue "E22" '  +0.64%  '

# Row 23
Set-TextValue "D23" '0.9989'
Set-TextValue "E23" '  -0.29%  '

# Row 24
Set-TextValue "D24" '7.755'
Set-TextValue "E24" '  +4.87%  '

# Row 25
Set-TextValue "D25" '0.1487'
Set-TextValue "E25" '  +3.16%  '

# Row 26
Set-TextValue "D26" '162.88'
Set-TextValue "E26" '  +1.03%  '

# Row 27
Set-TextValue "D27" '8.998'
Set-TextValue "E27" '  +0.53%  '

# Row 28
Set-TextValue "D28" '18.28'
Set-TextValue "E28" '  +1.30%  '

# Row 29
Set-TextValue "D29" '1.955'
Set-TextValue "E29" '  -0.34%  '

# Row 30
Set-TextValue "D30" '1.364'
Set-TextValue "E30" '  -4.76%  '

# Row 31
Set-TextValue "D31" '1.482'
Set-TextValue "E31" '  +0.11%  '

# Row 32
Set-TextValue "D32" '4.355'
Set-TextValue "E32" '  -0.69%  '

# Row 33
Set-TextValue "D33" '4.107'
Set-TextValue "E33" '  +1.46%  '

# Row 34
Set-TextValue "D34" '0.05248'
Set-TextValue "E34" '  +1.37%  '

# Row 35
Set-TextValue "D35" '1.197'
Set-TextValue "E35" '  +2.73%  '

# Row 36
Set-TextValue "D36" '0.7190'
Set-TextValue "E36" '  +1.97%  '

# Row 37
Set-TextValue "E37" '  +0.28%  '

# Row 38
Set-TextValue "E38" '  -0.08%  '

# Row 39
Set-TextValue "D39" '0.01861'
Set-TextValue "E39" '  +1.24%  '

# Row 40
Set-TextValue "D40" '2.711'
Set-TextValue "E40" '  -0.68%  '

# Row 41
Set-TextValue "D41" '1.179.26'
Set-TextValue "E41" '  +4.16%  '

# Row 42
Set-TextValue "D42" '0.9118'

# Row 43
Set-TextValue "B43" 'FraxShare'
Set-TextValue "C43" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D43" '5.992'
Set-TextValue "E43" '  +2.14%  '

# Row 44
Set-TextValue "B44" 'Aave'
Set-TextValue "C44" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D44" '72.03'
Set-TextValue "E44" '  +2.79%  '

# Row 45
Set-TextValue "D45" '0.4319'
Set-TextValue "E45" '  +1.31%  '

# Row 46
Set-TextValue "D46" '0.9987'
Set-TextValue "E46" '  -0.14%  '

# Row 47
Set-TextValue "D47" '102.47'
Set-TextValue "E47" '  +0.27%  '

# Row 48
Set-TextValue "D48" '0.5342'
Set-TextValue "E48" '  -1.74%  '

# Row 49
Set-TextValue "D49" '1.779'
Set-TextValue "E49" '  +1.04%  '

# Row 50
Set-TextValue "D50" '2.884'
Set-TextValue "E50" '  +5.20%  '

# Row 51
Set-TextValue "D51" '9.228'
Set-TextValue "E51" '  +0.69%  '
